{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"Echoes of Eternity\" -> \"The Profound Insights of\n// Mathematics and Its Widespread Impact\" rewrite described by the\n// commit diff.\n//\n// Strategy: for each old/new text pair, use body.search() to locate the\n// (unique) old text anywhere in the document body, then call\n// insertText(newText, \"Replace\") on the found range. This keeps the\n// surrounding runs / their formatting (rFonts/color/sz) untouched and\n// preserves straight apostrophes verbatim (no smart-quote autocorrect).\n\nconst body = context.document.body;\n\nasync function replaceText(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    console.log(\"NOT FOUND: \" + oldText);\n    return;\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1. Title\nawait replaceText(\n  \"Echoes of Eternity - A Celestial Symphony\",\n  \"The Profound Insights of Mathematics and Its Widespread Impact\"\n);\n\n// 2. Author name\nawait replaceText(\"Jaime Torres\", \"Dr. Albert Clayton\");\n\n// 3. Email local-part / domain / tld (already three separate runs in the doc)\nawait replaceText(\"jm\", \"ac\");\nawait replaceText(\"torres@umich\", \"claytonphd@protonmail\");\nawait replaceText(\"edu\", \"com\");\n\n// 4. Body paragraph 1 (first block, before the line breaks)\nawait replaceText(\n  \"Drawn by a cosmic magnetism, we are entranced by the silent yet resounding melodies of the universe\",\n  \"Mathematics, the language of the universe, unravels the intricate patterns and symmetries that permeate existence\"\n);\n\nawait replaceText(\n  \" The very fabric of our universe vibrates with energy, orchestrating a Symphony of the Cosmos\",\n  \" It is a boundless realm of exploration, where abstract concepts converge with tangible applications, illuminating the world around us\"\n);\n\nawait replaceText(\n  \" Throughout history, across cultures, humans have been bewitched by the celestial spectacle. From ancient astrologers to modern cosmologists, we attempt to decipher the harmonic frequencies, pulsating and interwoven, that compose this celestial ballet. Our telescopes and observatories act as divining rods, harmonically tuned to pick up these cosmic tunes. The history of space exploration can be seen as an effort to conduct this symphony, to understand and contribute to the grand cosmic opera, wherein celestial bodies play their individual musical parts\",\n  \" From the cosmos' vast expanses to the intricate designs of nature, mathematics provides a framework for comprehending and harnessing the universe's underlying forces\"\n);\n\n// 5. Body paragraph 1, second block (after first <br><br>)\nawait replaceText(\n  \"The celestial ambiance is woven with gravitational lullabies, where galaxies and stars dance in intricate waltzes\",\n  \"Immersed in a world governed by numbers and equations, we discern the rhythmic harmonies of mathematical principles echoing throughout our lives\"\n);\n\nawait replaceText(\n  \" The frequency of light from distant stars, the sonata of solar flares, the radioactive hum, and the melodic murmur of black holes - all these sounds, if transposed to human ears, would form a cacophony\",\n  \" The Pythagorean theorem unveils the beauty of geometric relationships, guiding architects and engineers in constructing awe-inspiring structures\"\n);\n\nawait replaceText(\n  \" Yet, this vibrant chaos has structure, order, and rhythm, waiting to be analyzed and decoded. Our scientific instruments become our musical instruments, lending us ears to unravel the cosmic score. As we probe deeper into the mysteries of the universe, we are becoming attuned to its musicality\",\n  \" Calculus, a symphony of change, empowers scientists to model complex phenomena, opening doors to novel technological advancements\"\n);\n\n// 6. Body paragraph 1, third block (after second <br><br>)\nawait replaceText(\n  \"On Earth, biology offers a terrestrial echo of the universe's Symphony\",\n  \"Mathematics serves as a venerable instrument of discovery, propelling humanity's quest for knowledge\"\n);\n\nawait replaceText(\n  \" Human cells exhibit rhythmic metabolic processes mimicking the pulsation of stars, and life's intricate system of interactions between organisms mirrors the celestial mechanics of planetary configurations\",\n  \" It unveils the secrets of the cosmos, unraveling the mysteries of celestial bodies and guiding astronauts through the vast expanse of space\"\n);\n\nawait replaceText(\n  \" By comprehending Earth's biological harmony, we deepen our understanding of the cosmic concord\",\n  \" It unlocks the enigmas of subatomic particles, empowering physicists to explore the fundamental building blocks of matter\"\n);\n\nawait replaceText(\n  \" Each species, like a unique instrument, adds its distinct timbre to the terrestrial symphony, contributing a unique rhythm to the symphony of life\",\n  \" Mathematics reveals the intricate machinery of life, enabling biologists to decipher the genetic code and unravel the complexities of the human body\"\n);\n\n// 7. Summary paragraph\nawait replaceText(\n  \"In essence, the universe reverberates with an unspoken melody, an intrinsic cohesion between celestial bodies, energy fields, and life itself\",\n  \"Mathematics, an intellectual odyssey, unveils the universe's mysteries, propels technological advancements, and underpins our understanding of the cosmos, nature, and life itself\"\n);\n\nawait replaceText(\n  \" Our goal is to capture this celestial symphony, to understand the rhythm of the universe, and to recognize our own place within its vast performance. By listening to the murmurs of cosmos, we come closer to comprehending the fabric of the universe and our own role within it\",\n  \" It remains an instrumental force in shaping our world, an enduring testament to the power of human intellect\"\n);\n\n// 8. Trailing empty paragraph added at the very end of the body (after the\n// Summary paragraph).\nbody.paragraphs.getLast().insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop script (PowerShell-style).\n# Applies the \"Echoes of Eternity\" -> \"The Profound Insights of Mathematics...\"\n# rewrite described by the commit diff.\n#\n# Strategy: for each old/new text pair, use Find to locate the (unique)\n# old text anywhere in the document body, then assign the new text to\n# the found Range's .Text property directly (NOT Find.Replacement.Text /\n# wdReplaceAll) so that Word's Find-and-Replace \"smart quote\" autocorrect\n# never silently mutates straight apostrophes into curly ones.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param(\n        [string]$Old,\n        [string]$New\n    )\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $Old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if ($found) {\n        $find.Parent.Text = $New\n    } else {\n        Write-Output \"NOT FOUND: $Old\"\n    }\n}\n\n# 1. Title\nReplace-Text \"Echoes of Eternity - A Celestial Symphony\" \"The Profound Insights of Mathematics and Its Widespread Impact\"\n\n# 2. Author name\nReplace-Text \"Jaime Torres\" \"Dr. Albert Clayton\"\n\n# 3. Email local-part / domain / tld (kept as three separate runs already in the doc)\nReplace-Text \"jm\" \"ac\"\nReplace-Text \"torres@umich\" \"claytonphd@protonmail\"\nReplace-Text \"edu\" \"com\"\n\n# 4. Body paragraph 1 (first block, before the line breaks)\nReplace-Text \"Drawn by a cosmic magnetism, we are entranced by the silent yet resounding melodies of the universe\" \"Mathematics, the language of the universe, unravels the intricate patterns and symmetries that permeate existence\"\n\nReplace-Text \" The very fabric of our universe vibrates with energy, orchestrating a Symphony of the Cosmos\" \" It is a boundless realm of exploration, where abstract concepts converge with tangible applications, illuminating the world around us\"\n\nReplace-Text \" Throughout history, across cultures, humans have been bewitched by the celestial spectacle. From ancient astrologers to modern cosmologists, we attempt to decipher the harmonic frequencies, pulsating and interwoven, that compose this celestial ballet. Our telescopes and observatories act as divining rods, harmonically tuned to pick up these cosmic tunes. The history of space exploration can be seen as an effort to conduct this symphony, to understand and contribute to the grand cosmic opera, wherein celestial bodies play their individual musical parts\" \" From the cosmos' vast expanses to the intricate designs of nature, mathematics provides a framework for comprehending and harnessing the universe's underlying forces\"\n\n# 5. Body paragraph 1, second block (after first <br><br>)\nReplace-Text \"The celestial ambiance is woven with gravitational lullabies, where galaxies and stars dance in intricate waltzes\" \"Immersed in a world governed by numbers and equations, we discern the rhythmic harmonies of mathematical principles echoing throughout our lives\"\n\nReplace-Text \" The frequency of light from distant stars, the sonata of solar flares, the radioactive hum, and the melodic murmur of black holes - all these sounds, if transposed to human ears, would form a cacophony\" \" The Pythagorean theorem unveils the beauty of geometric relationships, guiding architects and engineers in constructing awe-inspiring structures\"\n\nReplace-Text \" Yet, this vibrant chaos has structure, order, and rhythm, waiting to be analyzed and decoded. Our scientific instruments become our musical instruments, lending us ears to unravel the cosmic score. As we probe deeper into the mysteries of the universe, we are becoming attuned to its musicality\" \" Calculus, a symphony of change, empowers scientists to model complex phenomena, opening doors to novel technological advancements\"\n\n# 6. Body paragraph 1, third block (after second <br><br>)\nReplace-Text \"On Earth, biology offers a terrestrial echo of the universe's Symphony\" \"Mathematics serves as a venerable instrument of discovery, propelling humanity's quest for knowledge\"\n\nReplace-Text \" Human cells exhibit rhythmic metabolic processes mimicking the pulsation of stars, and life's intricate system of interactions between organisms mirrors the celestial mechanics of planetary configurations\" \" It unveils the secrets of the cosmos, unraveling the mysteries of celestial bodies and guiding astronauts through the vast expanse of space\"\n\nReplace-Text \" By comprehending Earth's biological harmony, we deepen our understanding of the cosmic concord\" \" It unlocks the enigmas of subatomic particles, empowering physicists to explore the fundamental building blocks of matter\"\n\nReplace-Text \" Each species, like a unique instrument, adds its distinct timbre to the terrestrial symphony, contributing a unique rhythm to the symphony of life\" \" Mathematics reveals the intricate machinery of life, enabling biologists to decipher the genetic code and unravel the complexities of the human body\"\n\n# 7. Summary paragraph\nReplace-Text \"In essence, the universe reverberates with an unspoken melody, an intrinsic cohesion between celestial bodies, energy fields, and life itself\" \"Mathematics, an intellectual odyssey, unveils the universe's mysteries, propels technological advancements, and underpins our understanding of the cosmos, nature, and life itself\"\n\nReplace-Text \" Our goal is to capture this celestial symphony, to understand the rhythm of the universe, and to recognize our own place within its vast performance. By listening to the murmurs of cosmos, we come closer to comprehending the fabric of the universe and our own role within it\" \" It remains an instrumental force in shaping our world, an enduring testament to the power of human intellect\"\n\n# 8. Trailing empty paragraph added at the very end of the body (after the Summary paragraph).\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n"}
